$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.890.75"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.638.43"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5096"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06440"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("E10").Value = "  +4.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.638.57"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "1.863.07"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5613"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "0.0₅7701"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "25.889.99"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.393"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.966"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.172"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.782"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1234"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.852"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04972"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.307"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.256"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.569"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.389"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.575"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.137.29"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.473"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8014"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  -3.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4267"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.788"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05072"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
